$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): extend with two new columns P and Q, matching the
# existing header style (same as O1 -> bold/bordered style "1").
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows (2-25): swap the I/K/M/O values (continuing the repeating
# 2,2,1 pattern) and append the two new data columns P and Q with value 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
